$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "Assasin" -> "Assassin"
$ws.Range("B2").Value = "Assassin"

# Add the new "Article" column (E) with the grammatical article for each character's name
$ws.Range("E1").Value = "Article"
$ws.Range("E2").Value = "L'"
$ws.Range("E3").Value = "Le "
$ws.Range("E4").Value = "Le "
$ws.Range("E5").Value = "Le "
$ws.Range("E6").Value = "L'"
$ws.Range("E7").Value = "Le "
$ws.Range("E8").Value = "L'"
$ws.Range("E9").Value = "La "

# Match the formatting of the rest of the header/data (left/top aligned, like column A)
$ws.Range("A1").Copy()
$ws.Range("E1:E9").PasteSpecial(-4122)

# Leave the selection where the author ended up after entering the new data
$ws.Range("E10").Select()
